$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.835.59"
$ws.Range("E2").Value = "  -0.87%  "

# Row 3
$ws.Range("D3").Value = "2.093.71"
$ws.Range("E3").Value = "  +2.12%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").Value = "245.04"
$ws.Range("E5").Value = "  -1.38%  "

# Row 6
$ws.Range("D6").Value = "0.656"
$ws.Range("E6").Value = "  -1.47%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "54.16"
$ws.Range("E8").Value = "  -5.15%  "

# Row 9
$ws.Range("D9").Value = "59.09"
$ws.Range("E9").Value = "  -1.63%  "

# Row 10
$ws.Range("D10").Value = "0.367"
$ws.Range("E10").Value = "  -4.60%  "

# Row 11
$ws.Range("D11").Value = "0.0767"
$ws.Range("E11").Value = "  -2.20%  "

# Row 12
$ws.Range("E12").Value = "  +1.08%  "

# Row 13
$ws.Range("E13").Value = "  +2.51%  "

# Row 14
$ws.Range("D14").Value = "14.96"
$ws.Range("E14").Value = "  -7.64%  "

# Row 15
$ws.Range("D15").Value = "2.396.83"
$ws.Range("E15").Value = "  +2.13%  "

# Row 16
$ws.Range("E16").Value = "  -4.69%  "

# Row 17
$ws.Range("D17").Value = "2.057.21"
$ws.Range("E17").Value = "  +0.32%  "

# Row 18
$ws.Range("D18").Value = "36.777.61"
$ws.Range("E18").Value = "  -1.05%  "

# Row 19
$ws.Range("D19").Value = "17.20"
$ws.Range("E19").Value = "  -8.94%  "

# Row 20
$ws.Range("D20").Value = "72.71"
$ws.Range("E20").Value = "  -2.55%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0882"
$ws.Range("E21").Value = "  -1.75%  "

# Row 22
$ws.Range("D22").Value = "5.46"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").Value = "240.00"
$ws.Range("E23").Value = "  +1.18%  "

# Row 24
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  -3.67%  "

# Row 26
$ws.Range("E26").Value = "  -0.13%  "

# Row 27
$ws.Range("D27").Value = "2.16"
$ws.Range("E27").Value = "  -0.86%  "

# Row 28
$ws.Range("D28").Value = "166.93"
$ws.Range("E28").Value = "  -1.91%  "

# Row 29
$ws.Range("D29").Value = "21.08"
$ws.Range("E29").Value = "  +4.35%  "

# Row 30
$ws.Range("E30").Value = "  -1.98%  "

# Row 31
$ws.Range("D31").Value = "5.22"
$ws.Range("E31").Value = "  +3.98%  "

# Row 32
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  -0.35%  "

# Row 33
$ws.Range("D33").Value = "4.73"
$ws.Range("E33").Value = "  +4.39%  "

# Row 34
$ws.Range("D34").Value = "0.0608"
$ws.Range("E34").Value = "  -2.58%  "

# Row 35
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").Value = "  +7.75%  "

# Row 36
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("E37").Value = "  +3.01%  "

# Row 38
$ws.Range("D38").Value = "0.0824"
$ws.Range("E38").Value = "  -6.72%  "

# Row 39
$ws.Range("D39").Value = "1.28"
$ws.Range("E39").Value = "  -5.16%  "

# Row 40
$ws.Range("E40").Value = "  +0.50%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0221"
$ws.Range("E41").Value = "  -1.39%  "

# Row 42
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").Value = "4.89"
$ws.Range("E42").Value = "  -8.27%  "

# Row 43
$ws.Range("D43").Value = "0.0963"
$ws.Range("E43").Value = "  -2.11%  "

# Row 44
$ws.Range("D44").Value = "96.60"
$ws.Range("E44").Value = "  +0.41%  "

# Row 45
$ws.Range("D45").Value = "2.86"
$ws.Range("E45").Value = "  -7.67%  "

# Row 46
$ws.Range("E46").Value = "  +14.39%  "

# Row 47
$ws.Range("D47").Value = "1.412.34"
$ws.Range("E47").Value = "  +10.85%  "

# Row 48
$ws.Range("D48").Value = "16.07"
$ws.Range("E48").Value = "  -8.71%  "

# Row 49
$ws.Range("D49").Value = "2.44"
$ws.Range("E49").Value = "  -0.61%  "

# Row 50
$ws.Range("D50").Value = "2.89"
$ws.Range("E50").Value = "  +1.56%  "

# Row 51
$ws.Range("D51").Value = "2.286.09"
$ws.Range("E51").Value = "  +2.33%  "
